$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# Insert a new (still empty) row before the current row 12; this shifts the
# existing rows 12..32 down to 13..33
$ws.Rows.Item(12).Insert()

# Grow the table to cover the inserted row (table now spans A1:C33)
$lo.Resize($ws.Range("A1:C33"))

# Append one more empty row at the bottom of the table
$newRow = $lo.ListRows.Add()
$lo.Resize($ws.Range("A1:C34"))

# Fill the appended row (row 34) first - matches the original authoring
# order, where this entry was typed in before the one inserted at row 12
$ws.Range("A34").Value = 14
$ws.Range("B34").Value = "RA14"
$ws.Range("C34").Value = "SD CS"

# Now fill the row that was inserted earlier (row 12) with the other mapping
$ws.Range("A12").Value = 14
$ws.Range("B12").Value = "RB14"
$ws.Range("C12").Value = "CS"

# Restore the active selection cell as left by the author
$ws.Range("B12").Select()
